$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Restyle A9:A12 (was "will not be implemented" style s=6) to the
#    "implemented" style (s=4), matching A8/A13's existing format.
# ---------------------------------------------------------------------------
$ws.Range("A8").Copy() | Out-Null
$ws.Range("A9:A12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# 2. Rows whose functions moved from "will not be implemented" to
#    " implemented oo and classic": restyle column A and column B to the
#    "implemented" look, and set column B's text accordingly.
# ---------------------------------------------------------------------------
$rowsNowImplemented = @(94, 100, 108, 112, 134, 157, 159, 200, 201)

foreach ($r in $rowsNowImplemented) {
    $ws.Range("A8").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

    $ws.Range("B8").Copy() | Out-Null
    $ws.Range("B$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
    $ws.Range("B$r").Value2 = " implemented oo and classic"
}

# Row 94 ends up with a plain (unstyled) B cell rather than the copied format.
$ws.Range("B94").ClearFormats() | Out-Null
$ws.Range("B94").Value2 = " implemented oo and classic"

# ---------------------------------------------------------------------------
# 3. Update the "Comments" column text for the rows above (C94 is set last
#    among the two new distinct strings so they are appended to the shared
#    string table in the same order as the reference edit).
# ---------------------------------------------------------------------------
$ws.Range("C159").Value2 = "Implemented indirectly in the equote() method and ooSQLiteEnquote()"
$ws.Range("C200").Value2 = "Implemented indirectly in the equote() method and ooSQLiteEnquote()"
$ws.Range("C201").Value2 = "Implemented indirectly in the equote() method and ooSQLiteEnquote()"

$ws.Range("C94").Value2  = "Implemented indirectly, used by implementation code"
$ws.Range("C100").Value2 = "Implemented indirectly, used by implementation code"
$ws.Range("C108").Value2 = "Implemented indirectly, used by implementation code"
$ws.Range("C134").Value2 = "Implemented indirectly, used by implementation code"
$ws.Range("C157").Value2 = "Implemented indirectly, used by implementation code"
# C112 keeps its original comment ("Of no use in Rexx") - no change needed.

# ---------------------------------------------------------------------------
# 4. Widen column C slightly and move the frozen-pane view / selection.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 66.140625

$ws.Application.ActiveWindow.ScrollRow = 65
$ws.Range("B65").Select() | Out-Null
